$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Turn the "Date:" / "Name:" paragraphs into a 2x2 table (each label in
#    the left column, an empty cell to the right) that keeps the original
#    paragraph formatting (BodyText style + run language) inside every cell.
# ---------------------------------------------------------------------------
$dateRange = $null
$nameRange = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Date:")) {
        $dateRange = $p.Range
    } elseif ($t.StartsWith("Name:") -and $nameRange -eq $null) {
        $nameRange = $p.Range
    }
    if ($dateRange -ne $null -and $nameRange -ne $null) {
        break
    }
}

$tableRange = $d.Range($dateRange.Start, $nameRange.End)

$tableXml = '<w:tbl>' +
    '<w:tblPr>' +
        '<w:tblStyle w:val="TableGrid"/>' +
        '<w:tblW w:w="9498" w:type="dxa"/>' +
        '<w:tblBorders>' +
            '<w:top w:val="none" w:sz="0" w:space="0" w:color="auto"/>' +
            '<w:left w:val="none" w:sz="0" w:space="0" w:color="auto"/>' +
            '<w:bottom w:val="none" w:sz="0" w:space="0" w:color="auto"/>' +
            '<w:right w:val="none" w:sz="0" w:space="0" w:color="auto"/>' +
            '<w:insideH w:val="none" w:sz="0" w:space="0" w:color="auto"/>' +
            '<w:insideV w:val="none" w:sz="0" w:space="0" w:color="auto"/>' +
        '</w:tblBorders>' +
        '<w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/>' +
    '</w:tblPr>' +
    '<w:tblGrid>' +
        '<w:gridCol w:w="877"/>' +
        '<w:gridCol w:w="8621"/>' +
    '</w:tblGrid>' +
    '<w:tr>' +
        '<w:tc>' +
            '<w:tcPr><w:tcW w:w="877" w:type="dxa"/></w:tcPr>' +
            '<w:p>' +
                '<w:pPr><w:pStyle w:val="BodyText"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
                '<w:r><w:rPr><w:lang w:val="en-GB"/></w:rPr><w:t>Date:</w:t></w:r>' +
            '</w:p>' +
        '</w:tc>' +
        '<w:tc>' +
            '<w:tcPr><w:tcW w:w="8621" w:type="dxa"/></w:tcPr>' +
            '<w:p>' +
                '<w:pPr><w:pStyle w:val="BodyText"/><w:rPr><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
            '</w:p>' +
        '</w:tc>' +
    '</w:tr>' +
    '<w:tr>' +
        '<w:tc>' +
            '<w:tcPr><w:tcW w:w="877" w:type="dxa"/></w:tcPr>' +
            '<w:p>' +
                '<w:pPr><w:pStyle w:val="BodyText"/><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr>' +
                '<w:r><w:rPr><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve">Name: </w:t></w:r>' +
            '</w:p>' +
        '</w:tc>' +
        '<w:tc>' +
            '<w:tcPr><w:tcW w:w="8621" w:type="dxa"/></w:tcPr>' +
            '<w:p>' +
                '<w:pPr><w:pStyle w:val="BodyText"/><w:rPr><w:lang w:val="nl-NL"/></w:rPr></w:pPr>' +
            '</w:p>' +
        '</w:tc>' +
    '</w:tr>' +
'</w:tbl>'

$tablePkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
    '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
            '<pkg:xmlData>' +
                '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                    '<w:body>' + $tableXml + '</w:body>' +
                '</w:document>' +
            '</pkg:xmlData>' +
        '</pkg:part>' +
    '</pkg:package>'

$tableRange.InsertXML($tablePkg)

# ---------------------------------------------------------------------------
# 2) Move <w:lastRenderedPageBreak/> from the "Trailing Annual Dividend
#    Yield" bullet to the "Profit Margin" bullet (it now sits right before
#    the "Profit Margin" run's text instead of the "Trailing Annual..." one).
# ---------------------------------------------------------------------------
$profitRange = $null
$trailingRange = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t.StartsWith("Profit Margin:")) {
        $profitRange = $p.Range
    } elseif ($t.StartsWith("Trailing Annual Dividend Yield:")) {
        $trailingRange = $p.Range
    }
}

$profitParaXml = '<w:p w14:paraId="5C281FE2" w14:textId="77777777" w:rsidR="00555B8C" w:rsidRPr="00555B8C" w:rsidRDefault="00555B8C" w:rsidP="00555B8C">' +
    '<w:pPr><w:pStyle w:val="BodyText"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="29"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="00555B8C"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-GB"/></w:rPr>' +
        '<w:lastRenderedPageBreak/>' +
        '<w:t>Profit Margin: This tells us how much profit a company is making from the money it earns. If the number is high, it means the company is making a good amount of profit.</w:t>' +
    '</w:r>' +
'</w:p>'

$trailingParaXml = '<w:p w14:paraId="0A5C0AE7" w14:textId="77777777" w:rsidR="00555B8C" w:rsidRPr="00555B8C" w:rsidRDefault="00555B8C" w:rsidP="00555B8C">' +
    '<w:pPr><w:pStyle w:val="BodyText"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="29"/></w:numPr><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-GB"/></w:rPr></w:pPr>' +
    '<w:r w:rsidRPr="00555B8C"><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-GB"/></w:rPr>' +
        '<w:t>Trailing Annual Dividend Yield: This tells us how much money a company pays to its shareholders as a percentage of the stock price. If the number is high, it means the company is giving more money back to its shareholders.</w:t>' +
    '</w:r>' +
'</w:p>'

function New-PartPkg([string]$bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
            '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
                '<pkg:xmlData>' +
                    '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                        '<w:body>' + $bodyXml + '</w:body>' +
                    '</w:document>' +
                '</pkg:xmlData>' +
            '</pkg:part>' +
        '</pkg:package>'
}

$trailingRange.InsertXML((New-PartPkg $trailingParaXml))
$profitRange.InsertXML((New-PartPkg $profitParaXml))

Write-Output "Edit complete"
